# Apply the 2023-07-03 cryptos data refresh (GitHub Actions scheduled update).
# Price (D) and Volume(1h) (E) values are refreshed for every coin row; a new
# coin (BitDAO) entered the top-50 at row 20, cascading several rows of
# Coin/Link (B/C) down, while the final row (previously Elrond) drops off.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'31.258.21"
$ws.Range("E2").Value = "'  +2.46%  "
# Row 3
$ws.Range("D3").Value = "'1.972.91"
$ws.Range("E3").Value = "'  +3.24%  "
# Row 4
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "'  +0.17%  "
# Row 5
$ws.Range("D5").Value = "'248.84"
$ws.Range("E5").Value = "'  +1.96%  "
# Row 6
$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "'  +0.23%  "
# Row 7
$ws.Range("D7").Value = "'0.4893"
$ws.Range("E7").Value = "'  +1.25%  "
# Row 8
$ws.Range("D8").Value = "'44.93"
$ws.Range("E8").Value = "'  +1.12%  "
# Row 9
$ws.Range("D9").Value = "'0.2960"
$ws.Range("E9").Value = "'  +2.55%  "
# Row 10
$ws.Range("D10").Value = "'0.06852"
$ws.Range("E10").Value = "'  +0.63%  "
# Row 11
$ws.Range("D11").Value = "'19.28"
$ws.Range("E11").Value = "'  -0.50%  "
# Row 12
$ws.Range("D12").Value = "'107.54"
$ws.Range("E12").Value = "'  -3.48%  "
# Row 13
$ws.Range("D13").Value = "'1.964.53"
$ws.Range("E13").Value = "'  +2.84%  "
# Row 14
$ws.Range("D14").Value = "'0.07791"
$ws.Range("E14").Value = "'  +3.04%  "
# Row 15
$ws.Range("D15").Value = "'5.465"
$ws.Range("E15").Value = "'  +1.67%  "
# Row 16
$ws.Range("D16").Value = "'0.7104"
$ws.Range("E16").Value = "'  +6.23%  "
# Row 17
$ws.Range("D17").Value = "'286.68"
$ws.Range("E17").Value = "'  -2.44%  "
# Row 18
$ws.Range("D18").Value = "'31.214.59"
$ws.Range("E18").Value = "'  +2.34%  "
# Row 19
$ws.Range("D19").Value = "'13.36"
$ws.Range("E19").Value = "'  +2.84%  "
# Row 20
$ws.Range("B20").Value = "'BitDAO"
$ws.Range("C20").Value = "'https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit"
$ws.Range("D20").Value = "'0.4977"
$ws.Range("E20").Value = "'  +12.85%  "
# Row 21
$ws.Range("B21").Value = "'ShibaInu"
$ws.Range("C21").Value = "'https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").Value = "'0.000007771"
$ws.Range("E21").Value = "'  +2.66%  "
# Row 22
$ws.Range("B22").Value = "'WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "'2.224.70"
$ws.Range("E22").Value = "'  +2.88%  "
# Row 23
$ws.Range("B23").Value = "'Uniswap"
$ws.Range("C23").Value = "'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").Value = "'5.646"
$ws.Range("E23").Value = "'  +2.69%  "
# Row 24
$ws.Range("B24").Value = "'Dai"
$ws.Range("C24").Value = "'https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").Value = "'1.002"
$ws.Range("E24").Value = "'  +0.24%  "
# Row 25
$ws.Range("B25").Value = "'BinanceUSD"
$ws.Range("C25").Value = "'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D25").Value = "'1.002"
$ws.Range("E25").Value = "'  +0.27%  "
# Row 26
$ws.Range("B26").Value = "'Chainlink"
$ws.Range("C26").Value = "'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D26").Value = "'6.658"
$ws.Range("E26").Value = "'  +3.96%  "
# Row 27
$ws.Range("B27").Value = "'Cosmos"
$ws.Range("C27").Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "'10.06"
$ws.Range("E27").Value = "'  +6.37%  "
# Row 28
$ws.Range("B28").Value = "'Monero"
$ws.Range("C28").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").Value = "'170.17"
$ws.Range("E28").Value = "'  +3.03%  "
# Row 29
$ws.Range("B29").Value = "'EthereumClassic"
$ws.Range("C29").Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "'20.17"
$ws.Range("E29").Value = "'  -0.68%  "
# Row 30
$ws.Range("B30").Value = "'LidoDAOToken"
$ws.Range("C30").Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D30").Value = "'2.203"
$ws.Range("E30").Value = "'  +6.13%  "
# Row 31
$ws.Range("B31").Value = "'Stellar"
$ws.Range("C31").Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").Value = "'0.1069"
$ws.Range("E31").Value = "'  +0.47%  "
# Row 32
$ws.Range("B32").Value = "'Toncoin"
$ws.Range("C32").Value = "'https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D32").Value = "'1.445"
$ws.Range("E32").Value = "'  +1.27%  "
# Row 33
$ws.Range("B33").Value = "'Filecoin"
$ws.Range("C33").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'4.859"
$ws.Range("E33").Value = "'  +19.89%  "
# Row 34
$ws.Range("B34").Value = "'InternetComputer(DFINITY)"
$ws.Range("C34").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").Value = "'4.531"
$ws.Range("E34").Value = "'  +9.66%  "
# Row 35
$ws.Range("B35").Value = "'Hedera"
$ws.Range("C35").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "'0.05088"
$ws.Range("E35").Value = "'  +2.28%  "
# Row 36
$ws.Range("B36").Value = "'ImmutableX"
$ws.Range("C36").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'0.7736"
$ws.Range("E36").Value = "'  +5.23%  "
# Row 37
$ws.Range("B37").Value = "'ARBITRUM"
$ws.Range("C37").Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").Value = "'1.175"
$ws.Range("E37").Value = "'  +3.76%  "
# Row 38
$ws.Range("D38").Value = "'0.02058"
$ws.Range("E38").Value = "'  +1.09%  "
# Row 39
$ws.Range("B39").Value = "'HuobiToken"
$ws.Range("C39").Value = "'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D39").Value = "'2.736"
$ws.Range("E39").Value = "'  +0.81%  "
# Row 40
$ws.Range("B40").Value = "'MXToken"
$ws.Range("C40").Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "'2.737"
$ws.Range("E40").Value = "'  +2.02%  "
# Row 41
$ws.Range("D41").Value = "'2.138"
$ws.Range("E41").Value = "'  +6.42%  "
# Row 42
$ws.Range("B42").Value = "'FraxShare"
$ws.Range("C42").Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "'6.455"
$ws.Range("E42").Value = "'  +11.83%  "
# Row 43
$ws.Range("B43").Value = "'Aave"
$ws.Range("C43").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").Value = "'74.34"
$ws.Range("E43").Value = "'  +7.57%  "
# Row 44
$ws.Range("B44").Value = "'TrustWalletToken"
$ws.Range("C44").Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").Value = "'0.8892"
$ws.Range("E44").Value = "'  +2.95%  "
# Row 45
$ws.Range("B45").Value = "'Quant"
$ws.Range("C45").Value = "'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").Value = "'110.39"
$ws.Range("E45").Value = "'  +1.13%  "
# Row 46
$ws.Range("B46").Value = "'TheSandbox"
$ws.Range("C46").Value = "'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D46").Value = "'0.4490"
$ws.Range("E46").Value = "'  +1.55%  "
# Row 47
$ws.Range("B47").Value = "'PaxDollar"
$ws.Range("C47").Value = "'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").Value = "'1.003"
$ws.Range("E47").Value = "'  +0.32%  "
# Row 48
$ws.Range("B48").Value = "'Aptos"
$ws.Range("C48").Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").Value = "'7.533"
$ws.Range("E48").Value = "'  +4.61%  "
# Row 49
$ws.Range("B49").Value = "'Maker"
$ws.Range("C49").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "'988.64"
$ws.Range("E49").Value = "'  +17.36%  "
# Row 50
$ws.Range("D50").Value = "'0.1274"
$ws.Range("E50").Value = "'  +3.85%  "
# Row 51
$ws.Range("B51").Value = "'EnergySwap"
$ws.Range("C51").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'9.428"
$ws.Range("E51").Value = "'  +2.81%  "

Write-Host "Applied 158 cell updates"
